# Auto-generated edit script: applies the value changes described by the
# Halicarnassus_Profits diff across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 79.5
$ws.Cells.Item(5, 10).Value = 56.5
$ws.Cells.Item(5, 12).Value = 56.5
$ws.Cells.Item(5, 14).Value = -286.5
$ws.Cells.Item(11, 8).Value = 26.6
$ws.Cells.Item(11, 9).Value = 26.6
$ws.Cells.Item(11, 11).Value = 26.6
$ws.Cells.Item(11, 13).Value = 113.4
$ws.Cells.Item(12, 8).Value = 167.125
$ws.Cells.Item(12, 9).Value = 165.66667
$ws.Cells.Item(12, 10).Value = 168
$ws.Cells.Item(12, 11).Value = 165.66667
$ws.Cells.Item(12, 12).Value = 168
$ws.Cells.Item(12, 13).Value = 4.333329999999989
$ws.Cells.Item(12, 14).Value = -508
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 5832
$ws.Cells.Item(40, 9).Value = 3749.625
$ws.Cells.Item(40, 11).Value = 3749.625
$ws.Cells.Item(40, 13).Value = -3574.625
$ws.Cells.Item(107, 8).Value = 111.666664
$ws.Cells.Item(107, 9).Value = 111.666664
$ws.Cells.Item(107, 11).Value = 111.666664
$ws.Cells.Item(107, 13).Value = 1808.333336
$ws.Cells.Item(111, 8).Value = 1243.8889
$ws.Cells.Item(111, 9).Value = 1086.875
$ws.Cells.Item(111, 10).Value = 2500
$ws.Cells.Item(111, 11).Value = 3260.625
$ws.Cells.Item(111, 12).Value = 7500
$ws.Cells.Item(111, 13).Value = -193.625
$ws.Cells.Item(111, 14).Value = -13634
$ws.Cells.Item(113, 8).Value = 1491.5
$ws.Cells.Item(113, 9).Value = 1489.8
$ws.Cells.Item(113, 11).Value = 1489.8
$ws.Cells.Item(113, 13).Value = 1764.2
$ws.Cells.Item(118, 8).Value = 1094.909
$ws.Cells.Item(118, 10).Value = 2999.6667
$ws.Cells.Item(118, 12).Value = 8999.000100000001
$ws.Cells.Item(118, 14).Value = -12313.0001
$ws.Cells.Item(131, 8).Value = 2021.7142
$ws.Cells.Item(131, 9).Value = 1963.75
$ws.Cells.Item(131, 11).Value = 5891.25
$ws.Cells.Item(131, 13).Value = -851.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1617.2174
$ws.Cells.Item(2, 9).Value = 1174
$ws.Cells.Item(2, 11).Value = 1174
$ws.Cells.Item(2, 13).Value = -1061
$ws.Cells.Item(32, 8).Value = 1727.5555
$ws.Cells.Item(32, 9).Value = 1634.0571
$ws.Cells.Item(32, 11).Value = 1634.0571
$ws.Cells.Item(32, 13).Value = -1347.0571
$ws.Cells.Item(61, 8).Value = 3792.2856
$ws.Cells.Item(61, 9).Value = 3090.2727
$ws.Cells.Item(61, 10).Value = 6366.3335
$ws.Cells.Item(61, 11).Value = 3090.2727
$ws.Cells.Item(61, 12).Value = 6366.3335
$ws.Cells.Item(61, 13).Value = -2878.2727
$ws.Cells.Item(61, 14).Value = -6790.3335
$ws.Cells.Item(80, 8).Value = 34073.332
$ws.Cells.Item(80, 9).Value = 2000
$ws.Cells.Item(80, 10).Value = 50110
$ws.Cells.Item(80, 11).Value = 2000
$ws.Cells.Item(80, 12).Value = 50110
$ws.Cells.Item(80, 14).Value = -52106
$ws.Cells.Item(80, 13).Value = -1002
$ws.Cells.Item(83, 8).Value = 34073.332
$ws.Cells.Item(83, 9).Value = 2000
$ws.Cells.Item(83, 10).Value = 50110
$ws.Cells.Item(83, 11).Value = 6000
$ws.Cells.Item(83, 12).Value = 150330
$ws.Cells.Item(83, 14).Value = -160314
$ws.Cells.Item(83, 13).Value = -1008
$ws.Cells.Item(116, 8).Value = 1617.2174
$ws.Cells.Item(116, 9).Value = 1174
$ws.Cells.Item(116, 11).Value = 1174
$ws.Cells.Item(116, 13).Value = 1120
$ws.Cells.Item(132, 8).Value = 2973.1
$ws.Cells.Item(132, 9).Value = 1390.2858
$ws.Cells.Item(132, 11).Value = 4170.857400000001
$ws.Cells.Item(132, 13).Value = -1640.857400000001
$ws.Cells.Item(136, 8).Value = 3792.2856
$ws.Cells.Item(136, 9).Value = 3090.2727
$ws.Cells.Item(136, 10).Value = 6366.3335
$ws.Cells.Item(136, 11).Value = 9270.8181
$ws.Cells.Item(136, 12).Value = 19099.0005
$ws.Cells.Item(136, 13).Value = -6720.8181
$ws.Cells.Item(136, 14).Value = -24199.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1617.2174
$ws.Cells.Item(3, 9).Value = 1174
$ws.Cells.Item(3, 11).Value = 1174
$ws.Cells.Item(3, 13).Value = -1060
$ws.Cells.Item(22, 8).Value = 369.6
$ws.Cells.Item(22, 10).Value = 814.3333
$ws.Cells.Item(22, 12).Value = 814.3333
$ws.Cells.Item(22, 14).Value = -1160.3333
$ws.Cells.Item(99, 8).Value = 3248.3
$ws.Cells.Item(99, 9).Value = 2717.375
$ws.Cells.Item(99, 11).Value = 2717.375
$ws.Cells.Item(99, 13).Value = -1219.375
$ws.Cells.Item(107, 8).Value = 3243.1785
$ws.Cells.Item(107, 9).Value = 1018.2941
$ws.Cells.Item(107, 10).Value = 6681.636
$ws.Cells.Item(107, 11).Value = 1018.2941
$ws.Cells.Item(107, 12).Value = 6681.636
$ws.Cells.Item(107, 13).Value = 901.7059
$ws.Cells.Item(107, 14).Value = -10521.636

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 322
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 322
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 322
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).Value = -546
$ws.Cells.Item(7, 8).Value = 70.117645
$ws.Cells.Item(7, 9).Value = 35.42857
$ws.Cells.Item(7, 11).Value = 35.42857
$ws.Cells.Item(7, 13).Value = 77.57142999999999
$ws.Cells.Item(22, 8).Value = 1981
$ws.Cells.Item(22, 9).Value = 395.5
$ws.Cells.Item(22, 11).Value = 395.5
$ws.Cells.Item(22, 13).Value = -45.5
$ws.Cells.Item(94, 8).Value = 4038.2727
$ws.Cells.Item(94, 9).Value = 1868.1666
$ws.Cells.Item(94, 10).Value = 6642.4
$ws.Cells.Item(94, 11).Value = 1868.1666
$ws.Cells.Item(94, 12).Value = 6642.4
$ws.Cells.Item(94, 13).Value = -1417.1666
$ws.Cells.Item(94, 14).Value = -7544.4
$ws.Cells.Item(99, 8).Value = 3141
$ws.Cells.Item(99, 9).Value = 2950.2727
$ws.Cells.Item(99, 11).Value = 2950.2727
$ws.Cells.Item(99, 13).Value = -1452.2727
$ws.Cells.Item(107, 8).Value = 558.6923
$ws.Cells.Item(107, 9).Value = 556.36365
$ws.Cells.Item(107, 11).Value = 556.36365
$ws.Cells.Item(107, 13).Value = 1363.63635
$ws.Cells.Item(126, 8).Value = 3141
$ws.Cells.Item(126, 9).Value = 2950.2727
$ws.Cells.Item(126, 11).Value = 8850.8181
$ws.Cells.Item(126, 13).Value = -6380.8181
$ws.Cells.Item(132, 8).Value = 2211.5334
$ws.Cells.Item(132, 9).Value = 1721
$ws.Cells.Item(132, 11).Value = 5163
$ws.Cells.Item(132, 13).Value = -2633

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 99.125
$ws.Cells.Item(10, 9).Value = 111.85714
$ws.Cells.Item(10, 11).Value = 335.57142
$ws.Cells.Item(10, 13).Value = -196.57142
$ws.Cells.Item(17, 8).Value = 1194.1428
$ws.Cells.Item(17, 9).Value = 393.16666
$ws.Cells.Item(17, 11).Value = 1179.49998
$ws.Cells.Item(17, 13).Value = -1010.49998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 474.4
$ws.Cells.Item(107, 9).Value = 474.4
$ws.Cells.Item(107, 11).Value = 474.4
$ws.Cells.Item(107, 13).Value = 1445.6
$ws.Cells.Item(126, 8).Value = 5192.875
$ws.Cells.Item(126, 9).Value = 5009.6
$ws.Cells.Item(126, 11).Value = 15028.8
$ws.Cells.Item(126, 13).Value = -12558.8
$ws.Cells.Item(132, 8).Value = 46688.848
$ws.Cells.Item(132, 9).Value = 65013.832
$ws.Cells.Item(132, 11).Value = 195041.496
$ws.Cells.Item(132, 13).Value = -192511.496

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 13).ClearContents()
$ws.Cells.Item(22, 8).Value = 1354.7693
$ws.Cells.Item(22, 9).Value = 516
$ws.Cells.Item(22, 11).Value = 516
$ws.Cells.Item(22, 13).Value = -221
$ws.Cells.Item(27, 8).Value = 1354.7693
$ws.Cells.Item(27, 9).Value = 516
$ws.Cells.Item(27, 11).Value = 516
$ws.Cells.Item(27, 13).Value = -409
$ws.Cells.Item(61, 8).Value = 3086
$ws.Cells.Item(61, 9).Value = 2259.0952
$ws.Cells.Item(61, 10).Value = 6559
$ws.Cells.Item(61, 11).Value = 2259.0952
$ws.Cells.Item(61, 12).Value = 6559
$ws.Cells.Item(61, 13).Value = -2057.0952
$ws.Cells.Item(61, 14).Value = -6963
$ws.Cells.Item(82, 8).Value = 5381.3335
$ws.Cells.Item(82, 10).Value = 5455.1665
$ws.Cells.Item(82, 12).Value = 5455.1665
$ws.Cells.Item(82, 14).Value = -6177.1665
$ws.Cells.Item(85, 8).Value = 5381.3335
$ws.Cells.Item(85, 10).Value = 5455.1665
$ws.Cells.Item(85, 12).Value = 5455.1665
$ws.Cells.Item(85, 14).Value = -7951.1665
$ws.Cells.Item(101, 8).Value = 5330.6665
$ws.Cells.Item(101, 10).Value = 5330.6665
$ws.Cells.Item(101, 12).Value = 5330.6665
$ws.Cells.Item(101, 14).Value = -11820.6665
$ws.Cells.Item(113, 8).Value = 3086
$ws.Cells.Item(113, 9).Value = 2259.0952
$ws.Cells.Item(113, 10).Value = 6559
$ws.Cells.Item(113, 11).Value = 2259.0952
$ws.Cells.Item(113, 12).Value = 6559
$ws.Cells.Item(113, 13).Value = -89.0952000000002
$ws.Cells.Item(113, 14).Value = -10899
$ws.Cells.Item(136, 8).Value = 3819.8333
$ws.Cells.Item(136, 9).Value = 3854.875
$ws.Cells.Item(136, 11).Value = 11564.625
$ws.Cells.Item(136, 13).Value = -9014.625
